$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalization fixes for particles (de/el/del -> De/El/Del) in municipality/state names
$ws.Range("B6").Value = "San Cristóbal De Las Casas"
$ws.Range("A12").Value = "Estado De México"
$ws.Range("B14").Value = "Apaseo El Alto"
$ws.Range("B18").Value = "Tulancingo De Bravo"
$ws.Range("B22").Value = "Unión De Tula"
$ws.Range("B34").Value = "Tepexi De Rodríguez"
$ws.Range("B41").Value = "Cosamaloapan De Carpio"
$ws.Range("B43").Value = "Paso Del Macho"

# Delete footer/metadata rows that are no longer needed.
# Delete from the bottom up so row numbers of earlier rows stay valid.
$ws.Range("A476:A480").EntireRow.Delete()
$ws.Range("A48:A52").EntireRow.Delete()

$ws.Range("A1:D46").Select()
